# Updated JPN model - 2025-09-04 00:57
# Apply edits to the "Misc" worksheet (a new wind-technology split row,
# a new "life" override row, and related label tweaks).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

# --- 1. Insert a blank row at 44 so the ~TFM_INS block (old rows 46-48)
#        shifts down to rows 47-49, making room for a new "E[_]WON*" row. ---
$ws.Rows("44:44").Insert()

# --- 2. Update the ~TFM_INS block (now rows 47-49) and append a new
#        "life" override row (50). ---
# F49 (was F48): add bioenergy to the technology list used by the "start" row.
$ws.Range("F49").Value = "solar,wind,coal,gas,nuclear,hydro,bioenergy"

# New row 50: a "-life" condition limiting LIFE overrides to
# coal/gas/nuclear/bioenergy processes.
$ws.Range("C50").Value = "life"
$ws.Range("D50").Value = 40
$ws.Range("E50").Value = "-life"
$ws.Range("F50").Value = "coal,gas,nuclear,bioenergy"

# E50 should keep the same "quote prefix" text style as E49 (so "-life"
# is stored as literal text, matching "-pasti" above it).
$ws.Range("E49").Copy()
$ws.Range("E50").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E50").Value = "-life"

# --- 3. Fix up the wind-offshore / wind-onshore rows in the ~TFM_TOPINS
#        block (rows 40-44). ---
# Row 43 used to be the combined "wind" row (E[_]W* / wind); split it into
# a dedicated offshore row ...
$ws.Range("C43").Value = "E[_]WOF*"
$ws.Range("D43").Value = "windoff"

# ... and add the new onshore row 44.
$ws.Range("C44").Value = "E[_]WON*"
$ws.Range("D44").Value = "windon"
$ws.Range("E44").Value = "IN"

# Row 41: ElcAgg_Wind now maps to the shorter ELC_wo* pattern (covers both
# onshore & offshore) instead of ELC_won*.
$ws.Range("D41").Value = "ELC_wo*"

# --- 4. Restore the view state (selection/scroll position) as closely as
#        possible. ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D42").Select()
